$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing quantities / prices (Charge Pumps) ---
$ws.Range("E3").Value2 = 10
$ws.Range("F3").Value2 = 0.041
$ws.Range("G3").Value2 = 0.41

$ws.Range("E9").Value2 = 20
$ws.Range("G9").Value2 = 1.6

# --- Add Total Price row ---
# Copy G12's format down into G13 first so the new total cell reuses the
# existing currency-style number format instead of minting a new style.
$ws.Range("G12").Copy($ws.Range("G13")) | Out-Null
$ws.Range("G13").Formula = "=SUM(G2:G12)"

$ws.Range("F13").Value2 = "Total Price"
$ws.Range("F13").Font.Bold = $true

# --- Update selection to match the author's final cursor position ---
$ws.Range("J4").Select() | Out-Null
